# Add a new worksheet "giordano_test" (liquid viscosity / density calc
# tests) right after Sheet1, populate it with the Giordano et al. test
# compositions, and update the selection on Sheet1 to the full data range.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# --- New worksheet, placed immediately after Sheet1 -----------------------
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet1)
$ws.Name = "giordano_test"

# --- Header row -------------------------------------------------------------
$headers = @("Label","SiO2","TiO2","Al2O3","FeO","Fe2O3","MgO","CaO","Na2O","K2O","P2O5","MnO","H2O","CO2","F","F2O")
for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# --- Data rows ----------------------------------------------------------------
$data = @(
    @("test_samp", 47.95, 1.67, 17.32, 10.24, 0.1, 5.76, 10.93, 3.45, 1.99, 0.51, 0.1, 2, 0.1, 0, 0),
    @("giordano_spreadsheet_default_comp", 62.4, 0.55000000000000004, 20.010000000000002, 0.03, 0, 3.22, 9.08, 3.52, 0.93, 0.12, 0.02, 2, 0, 0, 0.5),
    @("test_w_F", 47.95, 1.67, 17.32, 10.24, 0.1, 5.76, 10.93, 3.45, 1.99, 0.51, 0.1, 2, 0.1, 0.5, 0)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($col = 1; $col -le $row.Length; $col++) {
        $ws.Cells.Item($r + 2, $col).Value = $row[$col - 1]
    }
}

# --- Sheet1 selection now spans the whole used range, nothing "tabSelected" -
$sheet1.Range("A1:N2").Select() | Out-Null

# New sheet is the active one with the last-entered cell selected.
$ws.Activate() | Out-Null
$ws.Range("P4").Select() | Out-Null
